$wb = $excel.ActiveWorkbook

$wsPlanilha1 = $wb.Worksheets.Item("Planilha1")
$wsFeriados = $wb.Worksheets.Item("Feriados")

# Update C14:C25 values on Planilha1 from 0 to 2500000
$wsPlanilha1.Range("C14:C25").Value = 2500000

# Activate Planilha1 (first sheet) and set its selection
$wsPlanilha1.Activate()
$wsPlanilha1.Range("C13:C25").Select()

$wb.Save()
